$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts all existing data
# (and its formatting/column widths) from columns A:G to columns B:H.
$ws.Columns.Item(1).Insert()

# Header for the new "Day" column (bold, like the other headers).
$ws.Range("A1").Value = "Day"
$ws.Range("A1").Font.Bold = $true

# Day-number values: a literal seed value, then a relative formula series.
$ws.Range("A2").Value = 1
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4:A10").Formula = "=A3+1"

# Restore default (non-custom) row height for the header row.
$ws.Rows.Item(1).AutoFit()

# Clear the old cell selection so it reverts to the sheet default.
$ws.Range("A1").Select() | Out-Null

# Update the defined name that anchored the old Date header cell so it
# continues to point at the same logical cell, now shifted to column B.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_Hlk35725696") {
        $n.RefersTo = "=TompkinsCountyHealthDepartmentC!`$B`$1"
    }
}

# Recreate the sort state so its persisted range references follow the
# data that moved from columns A:G to columns B:H.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B10")) | Out-Null
$ws.Sort.SetRange($ws.Range("B2:H10"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

$wb.Save()
